$wb = $excel.ActiveWorkbook

# ProductLoanInput sheet: update the "shortname" value (B2) from the
# numeric 3524 to the text "353a".
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsInput.Range("B2").Value = "353a"

# Make ProductLoanInput the active/selected sheet (was ProductLoanOutput),
# with B7 as the selected cell.
$wsInput.Select()
$wsInput.Range("B7").Select()
